$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V1_AxonStim")

# Add two new rows of cell-library data (rows 40 and 41)
$ws.Range("A40").Value = "EB_042114_A"
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "AL"
$ws.Range("D40").Value = "PY"
$ws.Range("E40").Value = "2/3"
$ws.Range("F40").Value = 1
$ws.Rows.Item(40).RowHeight = 26.1

$ws.Range("A41").Value = "EB_042114_A"
$ws.Range("B41").Value = 2
$ws.Range("C41").Value = "PM"
$ws.Range("D41").Value = "PY"
$ws.Range("E41").Value = "2/3"
$ws.Range("F41").Value = 1
$ws.Rows.Item(41).RowHeight = 26.1

# Update the frozen pane view and active selection to reflect scrolling to the new rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("F41").Select()
